$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3 (pushes existing row 3 "country_of_residence" down to row 4, etc.)
$ws.Rows(3).Insert()

# Populate the new row 3 with the "project_name" field definition
$ws.Range("A3").Value = ""
$ws.Range("B3").Value = "project_name"
$ws.Range("C3").Value = "text"
$ws.Range("D3").Value = "Yes"
$ws.Range("E3").Value = ""

# Match styling: A3/B3/C3/D3/E3 pick up styles similar to neighboring body rows
$ws.Range("A3").Style = $ws.Range("A4").Style
$ws.Range("B3").Style = $ws.Range("B1").Style
$ws.Range("C3").Style = $ws.Range("B1").Style
$ws.Range("D3").Style = $ws.Range("B1").Style
$ws.Range("E3").Style = $ws.Range("E4").Style

# Move the selection to D3 as in the saved file, and reset the top-left scroll cell
$ws.Application.Goto($ws.Range("D3"))
$ws.Range("D3").Select()

Write-Output "done"
